$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 289.53845
$ws.Range("I12").Value = 336.4
$ws.Range("J12").Value = 133.33333
$ws.Range("K12").Value = 336.4
$ws.Range("L12").Value = 133.33333
$ws.Range("M12").Value = -166.4
$ws.Range("N12").Value = -473.33333
$ws.Range("H43").Value = 3081303
$ws.Range("I43").Value = 5131505
$ws.Range("J43").Value = 6000
$ws.Range("K43").Value = 5131505
$ws.Range("L43").Value = 6000
$ws.Range("M43").Value = -5131436
$ws.Range("N43").Value = -6138
$ws.Range("H69").Value = 9999
$ws.Range("J69").Value = 9999
$ws.Range("L69").Value = 29997
$ws.Range("N69").Value = -31745
$ws.Range("H72").Value = 9999
$ws.Range("J72").Value = 9999
$ws.Range("L72").Value = 89991
$ws.Range("N72").Value = -98727
$ws.Range("H80").Value = 704.6111
$ws.Range("I80").Value = 505.25
$ws.Range("J80").Value = 1103.3334
$ws.Range("K80").Value = 1515.75
$ws.Range("L80").Value = 3310.0002
$ws.Range("M80").Value = -517.75
$ws.Range("N80").Value = -5306.0002
$ws.Range("H83").Value = 704.6111
$ws.Range("I83").Value = 505.25
$ws.Range("J83").Value = 1103.3334
$ws.Range("K83").Value = 4547.25
$ws.Range("L83").Value = 9930.000599999999
$ws.Range("M83").Value = 444.75
$ws.Range("N83").Value = -19914.0006
$ws.Range("H116").Value = 25781174
$ws.Range("I116").Value = 38028228
$ws.Range("K116").Value = 38028228
$ws.Range("M116").Value = -38024786
$ws.Range("H125").Value = 4077.4
$ws.Range("J125").Value = 4077.4
$ws.Range("L125").Value = 36696.6
$ws.Range("N125").Value = -41616.6

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3508.2307
$ws.Range("I32").Value = 1804.4
$ws.Range("K32").Value = 1804.4
$ws.Range("M32").Value = -1517.4
$ws.Range("J63").Value = 4500
$ws.Range("L63").Value = 4500
$ws.Range("N63").Value = -5872
$ws.Range("J66").Value = 4500
$ws.Range("L66").Value = 22500
$ws.Range("N66").Value = -29364
$ws.Range("H74").Value = 10418448
$ws.Range("I74").Value = 20834582
$ws.Range("K74").Value = 20834582
$ws.Range("M74").Value = -20833708
$ws.Range("H77").Value = 10418448
$ws.Range("I77").Value = 20834582
$ws.Range("K77").Value = 104172910
$ws.Range("M77").Value = -104168542
$ws.Range("H132").Value = 30821.863
$ws.Range("I132").Value = 39286.066
$ws.Range("J132").Value = 12684.286
$ws.Range("K132").Value = 117858.198
$ws.Range("L132").Value = 38052.858
$ws.Range("M132").Value = -115328.198
$ws.Range("N132").Value = -43112.858

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1310.8889
$ws.Range("I86").Value = 1343.3846
$ws.Range("J86").Value = 1226.4
$ws.Range("K86").Value = 1343.3846
$ws.Range("L86").Value = 1226.4
$ws.Range("M86").Value = -220.3846000000001
$ws.Range("N86").Value = -3472.4
$ws.Range("H89").Value = 1310.8889
$ws.Range("I89").Value = 1343.3846
$ws.Range("J89").Value = 1226.4
$ws.Range("K89").Value = 6716.923000000001
$ws.Range("L89").Value = 6132
$ws.Range("M89").Value = -1100.923000000001
$ws.Range("N89").Value = -17364
$ws.Range("H99").Value = 1390449.1
$ws.Range("I99").Value = 2605432.5
$ws.Range("K99").Value = 2605432.5
$ws.Range("M99").Value = -2603934.5
$ws.Range("H105").Value = 2894.111
$ws.Range("I105").Value = 2584.3572
$ws.Range("K105").Value = 2584.3572
$ws.Range("M105").Value = -837.3571999999999
$ws.Range("H134").Value = 2637.4827
$ws.Range("I134").Value = 1591.1666
$ws.Range("K134").Value = 4773.4998
$ws.Range("M134").Value = -2238.4998

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 13891478
$ws.Range("I132").Value = 18520788
$ws.Range("K132").Value = 55562364
$ws.Range("M132").Value = -55559834

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 516.25
$ws.Range("J9").Value = 474.5
$ws.Range("L9").Value = 1423.5
$ws.Range("N9").Value = -1871.5
$ws.Range("H86").Value = 214.2
$ws.Range("I86").Value = 195
$ws.Range("J86").Value = 227
$ws.Range("K86").Value = 585
$ws.Range("L86").Value = 681
$ws.Range("M86").Value = 601
$ws.Range("N86").Value = -3053
$ws.Range("H89").Value = 214.2
$ws.Range("I89").Value = 195
$ws.Range("J89").Value = 227
$ws.Range("K89").Value = 1755
$ws.Range("L89").Value = 2043
$ws.Range("M89").Value = 4173
$ws.Range("N89").Value = -13899
$ws.Range("H113").Value = 936.6
$ws.Range("J113").Value = 971.25
$ws.Range("L113").Value = 2913.75
$ws.Range("N113").Value = -7253.75
$ws.Range("H121").Value = 501071.66
$ws.Range("J121").Value = 750808
$ws.Range("L121").Value = 2252424
$ws.Range("N121").Value = -2255044
$ws.Range("H122").Value = 724.125
$ws.Range("J122").Value = 699
$ws.Range("L122").Value = 6291
$ws.Range("N122").Value = -11191
$ws.Range("H129").Value = 1596.2
$ws.Range("J129").Value = 1515.8334
$ws.Range("L129").Value = 4547.5002
$ws.Range("N129").Value = -14547.5002
$ws.Range("H131").Value = 10206054
$ws.Range("J131").Value = 7577873
$ws.Range("L131").Value = 22733619
$ws.Range("N131").Value = -22743699
$ws.Range("H140").Value = 5765.625
$ws.Range("I140").Value = 6605.4
$ws.Range("K140").Value = 19816.2
$ws.Range("M140").Value = -14636.2

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1216523
$ws.Range("I80").Value = 3130159.5
$ws.Range("J80").Value = 20500.125
$ws.Range("K80").Value = 3130159.5
$ws.Range("L80").Value = 20500.125
$ws.Range("M80").Value = -3129161.5
$ws.Range("N80").Value = -22496.125
$ws.Range("H83").Value = 1216523
$ws.Range("I83").Value = 3130159.5
$ws.Range("J83").Value = 20500.125
$ws.Range("K83").Value = 15650797.5
$ws.Range("L83").Value = 102500.625
$ws.Range("M83").Value = -15645805.5
$ws.Range("N83").Value = -112484.625
$ws.Range("H102").Value = 20839768
$ws.Range("I102").Value = 27785552
$ws.Range("J102").Value = 2415.8333
$ws.Range("K102").Value = 27785552
$ws.Range("L102").Value = 2415.8333
$ws.Range("M102").Value = -27783930
$ws.Range("N102").Value = -5659.8333
$ws.Range("H122").Value = 2705.2856
$ws.Range("I122").Value = 1999.2
$ws.Range("K122").Value = 5997.6
$ws.Range("M122").Value = -3547.6
$ws.Range("H132").Value = 8107.364
$ws.Range("I132").Value = 7451.0713
$ws.Range("K132").Value = 22353.2139
$ws.Range("M132").Value = -19823.2139

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5610.3335
$ws.Range("I7").Value = 4641.143
$ws.Range("J7").Value = 9002.5
$ws.Range("K7").Value = 4641.143
$ws.Range("L7").Value = 9002.5
$ws.Range("M7").Value = -4529.143
$ws.Range("N7").Value = -9226.5
$ws.Range("H82").Value = 6252200
$ws.Range("I82").Value = 7814875
$ws.Range("K82").Value = 7814875
$ws.Range("M82").Value = -7814514
$ws.Range("H85").Value = 6252200
$ws.Range("I85").Value = 7814875
$ws.Range("K85").Value = 7814875
$ws.Range("M85").Value = -7813627
$ws.Range("H93").Value = 4041.125
$ws.Range("I93").Value = 4332.25
$ws.Range("K93").Value = 4332.25
$ws.Range("M93").Value = -3084.25
$ws.Range("H126").Value = 5610.3335
$ws.Range("I126").Value = 4641.143
$ws.Range("J126").Value = 9002.5
$ws.Range("K126").Value = 13923.429
$ws.Range("L126").Value = 27007.5
$ws.Range("M126").Value = -11453.429
$ws.Range("N126").Value = -31947.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3681.6052
$ws.Range("I122").Value = 3564.375
$ws.Range("K122").Value = 10693.125
$ws.Range("M122").Value = -8243.125
